# The deck's theme (ppt/theme/theme1.xml, used by the slide master and
# therefore by every slide) is switched from the "Integral" color scheme
# to the stock "Office Theme" color scheme. The font scheme and format
# scheme are already identical between the two themes in this deck, so
# only the twelve theme colors need to change.
#
# PowerPoint's Theme object doesn't expose a "swap in another theme part"
# verb directly, but it does expose the live color slots via
# ThemeColorScheme.Item(n).RGB (VBA/COM RGB ordering: R | (G<<8) | (B<<16)),
# in the fixed order dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink -
# exactly the <a:clrScheme> child order. Driving each slot to the "Office"
# palette reproduces the target clrScheme content.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# Office Theme palette (target), in clrScheme order, encoded as VBA RGB()
# (R + G*256 + B*65536):
#   1  dk1      000000 -> 0
#   2  lt1      FFFFFF -> 16777215
#   3  dk2      44546A -> 6968388
#   4  lt2      E7E6E6 -> 15132391
#   5  accent1  5B9BD5 -> 13998939
#   6  accent2  ED7D31 -> 3243501
#   7  accent3  A5A5A5 -> 10855845
#   8  accent4  FFC000 -> 49407
#   9  accent5  4472C4 -> 12874308
#   10 accent6  70AD47 -> 4697456
#   11 hlink    0563C1 -> 12673797
#   12 folHlink 954F72 -> 7491477
$officeRgb = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le $colors.Count; $i++) {
    $colors.Item($i).RGB = $officeRgb[$i - 1]
}
